$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (existing B/C shift to C/D)
$ws.Columns.Item(2).Insert()

# Populate the new column B with the StatQuery header/value
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.race IN ['BLACK_OR_AFRICAN_AMERICAN'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match the wrap-text style used by A2 for the new B2 cell
$ws.Range("B2").WrapText = $true

# Column B needs a width matching column A's (75.81640625 chars); columns A, C, D
# already retain their original widths after the column insert.
$ws.Columns.Item(2).ColumnWidth = 75

# Update selection/active cell
[void]$ws.Range("A2").Select()
